$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("expenses")

# New row of data (row 2). The values "78000" and "15200" must be stored
# as text (shared strings), not numbers, matching the source file.
# Temporarily mark the cells as Text before assigning so the numeric-
# looking strings aren't auto-converted to numbers, then clear the
# formatting back off so no extra style is left applied to the cells.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "78000"
$ws.Range("A2").ClearFormats()

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "15200"
$ws.Range("C2").ClearFormats()
